$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> Utenti)
$ws.Name = "Utenti"

# Explicit "auto-adjusted" column widths (stored OOXML "width" = ColumnWidth +
# 5/6 padding that Excel adds on top for the default Calibri 11 font, so the
# ColumnWidth values below are chosen to land exactly on the target widths).
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 22.166666666666668
$ws.Columns.Item(5).ColumnWidth = 17.166666666666668
$ws.Columns.Item(6).ColumnWidth = 30.166666666666668
$ws.Columns.Item(7).ColumnWidth = 20.166666666666668
$ws.Columns.Item(8).ColumnWidth = 68.16666666666667

# Refresh the user records (rows 2-11)
# Row 2
$ws.Range("A2").Value = "Daniele"
$ws.Range("B2").Value = "Cesaroni"
$ws.Range("C2").Value = "24/12/1993"
$ws.Range("D2").Value = "Vallonga"
$ws.Range("E2").Value = "BRTSTN13H48I163K"
$ws.Range("F2").Value = "giampaolomajewski@example.net"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "033241376"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = "Vicolo Costantino, 227`n95124, Catania (CT)"
$ws.Rows.Item(2).AutoFit()

# Row 3
$ws.Range("A3").Value = "Adelmo"
$ws.Range("B3").Value = "Greco"
$ws.Range("C3").Value = "16/10/2006"
$ws.Range("D3").Value = "Antonimina"
$ws.Range("E3").Value = "VRGLSN74P28H569O"
$ws.Range("F3").Value = "dinasemitecolo@example.com"
$ws.Range("G3").Value = "+39 052104004"
$ws.Range("H3").Value = "Stretto Argentero, 4`n86100, Campobasso (CB)"
$ws.Rows.Item(3).AutoFit()

# Row 4
$ws.Range("A4").Value = "Sole"
$ws.Range("B4").Value = "Marenzio"
$ws.Range("C4").Value = "13/07/1996"
$ws.Range("D4").Value = "Camigliano"
$ws.Range("E4").Value = "FSCGCN56B49C448I"
$ws.Range("F4").Value = "ysatriani@example.org"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "05857520547"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = "Vicolo Annamaria, 51`n98168, Pace (ME)"
$ws.Rows.Item(4).AutoFit()

# Row 5
$ws.Range("A5").Value = "Giuseppina"
$ws.Range("B5").Value = "Goldoni"
$ws.Range("C5").Value = "20/09/2004"
$ws.Range("D5").Value = "Morrovalle"
$ws.Range("E5").Value = "RGGLSS82R70L315J"
$ws.Range("F5").Value = "vsibilia@example.com"
$ws.Range("G5").Value = "+39 0344405064"
$ws.Range("H5").Value = "Vicolo Ceri, 61`n98062, Ficarra (ME)"
$ws.Rows.Item(5).AutoFit()

# Row 6
$ws.Range("A6").Value = "Annibale"
$ws.Range("B6").Value = "Tozzi"
$ws.Range("C6").Value = "21/08/1997"
$ws.Range("D6").Value = "Braone"
$ws.Range("E6").Value = "NTLRFL15D09I030A"
$ws.Range("F6").Value = "pellicomarisa@example.com"
$ws.Range("G6").Value = "+39 016573231"
$ws.Range("H6").Value = "Stretto Regge, 51`n47824, Torriana (RN)"
$ws.Rows.Item(6).AutoFit()

# Row 7
$ws.Range("A7").Value = "Mattia"
$ws.Range("B7").Value = "Gargallo"
$ws.Range("C7").Value = "25/03/2002"
$ws.Range("D7").Value = "Gualdo"
$ws.Range("E7").Value = "ZCCFRN86P61C746Y"
$ws.Range("F7").Value = "ngrassi@example.net"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "3625387919"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = "Vicolo Iolanda, 90 Piano 3`n20056, Trezzo Sull'Adda (MI)"
$ws.Rows.Item(7).AutoFit()

# Row 8
$ws.Range("A8").Value = "Paolo"
$ws.Range("B8").Value = "Bragadin"
$ws.Range("C8").Value = "21/07/1993"
$ws.Range("D8").Value = "Popiglio"
$ws.Range("E8").Value = "ZTTRTR56M50F256V"
$ws.Range("F8").Value = "pcarli@example.org"
$ws.Range("G8").Value = "+39 05562205505"
$ws.Range("H8").Value = "Rotonda Eva, 36 Appartamento 56`n28070, Nibbiola (NO)"
$ws.Rows.Item(8).AutoFit()

# Row 9
$ws.Range("A9").Value = "Vittoria"
$ws.Range("B9").Value = "Bossi"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "12/08/1989"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "Petrella Liri"
$ws.Range("E9").Value = "TRNRLD40M53A766W"
$ws.Range("F9").Value = "vittoriatolentino@example.com"
$ws.Range("G9").Value = "+39 3770522409"
$ws.Range("H9").Value = "Borgo Brambilla, 72 Appartamento 98`n47833, Morciano Di Romagna (RN)"
$ws.Rows.Item(9).AutoFit()

# Row 10
$ws.Range("A10").Value = "Luchino"
$ws.Range("B10").Value = "Bellini"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "12/03/1991"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "Giugliano In Campania"
$ws.Range("E10").Value = "NTTLSE00A13B429T"
$ws.Range("F10").Value = "enricosalvo@example.net"
$ws.Range("G10").Value = "+39 37769446312"
$ws.Range("H10").Value = "Via Zoppetti, 33 Piano 1`n54021, Treschietto (MS)"
$ws.Rows.Item(10).AutoFit()

# Row 11
$ws.Range("A11").Value = "Calogero"
$ws.Range("B11").Value = "Germano"
$ws.Range("C11").Value = "20/03/2000"
$ws.Range("D11").Value = "Cavoli"
$ws.Range("E11").Value = "TRTLSS64A22L730P"
$ws.Range("F11").Value = "ubaldo78@example.net"
$ws.Range("G11").Value = "+39 0353917638"
$ws.Range("H11").Value = "Viale Borrani, 66 Appartamento 16`n82036, Solopaca (BN)"
$ws.Rows.Item(11).AutoFit()
